$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K4").Value = 12
$ws.Range("L4").Value = 19.6

$ws.Range("K5").Value = 25.11
$ws.Range("L5").Value = 36.67

$ws.Range("K6").Value = 28
$ws.Range("L6").Value = 73.81999999999999

$ws.Range("K7").Value = 131.163
$ws.Range("L7").Value = 99.02

$ws.Range("K8").Value = 201.273
$ws.Range("L8").Value = 240.12
